# Added 4K timelapse toggle
# Default is set to 1080p, regardless of device -> mark the "OK" (blue) status
# for all device columns (B:F) on the "4k Time Lapse" row (row 36).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in "OK" for the 4k Time Lapse row across all device columns (B-F),
# matching the formatting used by the other rows in this table (blue text).
$rng = $ws.Range("B36:F36")
$rng.Value = "OK"
$rng.Font.Color = 12611584

# The row grows slightly shorter to match the compact row height used
# elsewhere in the sheet once it gets populated.
$ws.Rows.Item(36).RowHeight = 12.8

# Move the active selection, as left by the author after the edit.
[void]$ws.Range("H33").Select()
